$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 13346.667
$ws.Range("I6").Value = 13346.667
$ws.Range("K6").Value = 40040.001
$ws.Range("M6").Value = -39928.001
$ws.Range("H7").Value = 11249.25
$ws.Range("J7").Value = 12998.5
$ws.Range("L7").Value = 12998.5
$ws.Range("N7").Value = -13222.5
$ws.Range("H14").Value = 11249.25
$ws.Range("J14").Value = 12998.5
$ws.Range("L14").Value = 12998.5
$ws.Range("N14").Value = -13380.5
$ws.Range("H92").Value = 1526.7858
$ws.Range("I92").Value = 1123.7727
$ws.Range("J92").Value = 3004.5
$ws.Range("K92").Value = 1123.7727
$ws.Range("L92").Value = 3004.5
$ws.Range("M92").Value = 124.2273
$ws.Range("N92").Value = -5500.5
$ws.Range("H113").Value = 6263.923
$ws.Range("I113").Value = 6868.3
$ws.Range("K113").Value = 6868.3
$ws.Range("M113").Value = -3614.3
$ws.Range("H127").Value = 45469984
$ws.Range("I127").Value = 66688012
$ws.Range("K127").Value = 200064036
$ws.Range("M127").Value = -200059076
$ws.Range("H131").Value = 131889
$ws.Range("I131").Value = 5618.8887
$ws.Range("K131").Value = 16856.6661
$ws.Range("M131").Value = -11816.6661
$ws.Range("H135").Value = 2117.3872
$ws.Range("I135").Value = 2169.4348
$ws.Range("K135").Value = 19524.9132
$ws.Range("M135").Value = -16989.9132
$ws.Range("H138").Value = 6767.7
$ws.Range("J138").Value = 6936.7393
$ws.Range("L138").Value = 20810.2179
$ws.Range("N138").Value = -31090.2179
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41570.395
$ws.Range("I32").Value = 34535.094
$ws.Range("J32").Value = 68432.45
$ws.Range("K32").Value = 34535.094
$ws.Range("L32").Value = 68432.45
$ws.Range("M32").Value = -34248.094
$ws.Range("N32").Value = -69006.45
$ws.Range("H61").Value = 7827.625
$ws.Range("J61").Value = 8765.5
$ws.Range("L61").Value = 8765.5
$ws.Range("N61").Value = -9189.5
$ws.Range("H74").Value = 8201.368
$ws.Range("I74").Value = 8705.441000000001
$ws.Range("K74").Value = 8705.441000000001
$ws.Range("M74").Value = -7831.441000000001
$ws.Range("H77").Value = 8201.368
$ws.Range("I77").Value = 8705.441000000001
$ws.Range("K77").Value = 43527.205
$ws.Range("M77").Value = -39159.205
$ws.Range("H110").Value = 9721.904
$ws.Range("I110").Value = 5951.0586
$ws.Range("K110").Value = 5951.0586
$ws.Range("M110").Value = -3906.0586
$ws.Range("H132").Value = 6391.731
$ws.Range("I132").Value = 6147
$ws.Range("K132").Value = 18441
$ws.Range("M132").Value = -15911
$ws.Range("H136").Value = 7827.625
$ws.Range("J136").Value = 8765.5
$ws.Range("L136").Value = 26296.5
$ws.Range("N136").Value = -31396.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55561944
$ws.Range("I31").Value = 3940.125
$ws.Range("K31").Value = 3940.125
$ws.Range("M31").Value = -3645.125
$ws.Range("H34").Value = 55561944
$ws.Range("I34").Value = 3940.125
$ws.Range("K34").Value = 3940.125
$ws.Range("M34").Value = -3738.125
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H132").Value = 4096.3022
$ws.Range("I132").Value = 2859.6333
$ws.Range("K132").Value = 8578.8999
$ws.Range("M132").Value = -6048.8999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5555772
$ws.Range("I2").Value = 250.88889
$ws.Range("K2").Value = 1505.33334
$ws.Range("M2").Value = -1392.33334
$ws.Range("H7").Value = 168563.5
$ws.Range("J7").Value = 3496.6667
$ws.Range("L7").Value = 10490.0001
$ws.Range("N7").Value = -10714.0001
$ws.Range("H23").Value = 183.44444
$ws.Range("J23").Value = 150.5
$ws.Range("L23").Value = 451.5
$ws.Range("N23").Value = -921.5
$ws.Range("H34").Value = 189
$ws.Range("I34").Value = 189
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 567
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -483
$ws.Range("N34").ClearContents()
$ws.Range("H80").Value = 499.5
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 499.5
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H92").Value = 856
$ws.Range("I92").Value = 808
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2424
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -1176
$ws.Range("N92").Value = -5496
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10332
$ws.Range("H122").Value = 2042.4242
$ws.Range("I122").Value = 1548.1786
$ws.Range("K122").Value = 4644.5358
$ws.Range("M122").Value = -2194.5358
$ws.Range("H132").Value = 18716.328
$ws.Range("I132").Value = 23395.715
$ws.Range("J132").Value = 3430.3333
$ws.Range("K132").Value = 70187.145
$ws.Range("L132").Value = 10290.9999
$ws.Range("M132").Value = -67657.145
$ws.Range("N132").Value = -15350.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 920707.25
$ws.Range("I40").Value = 920707.25
$ws.Range("K40").Value = 920707.25
$ws.Range("M40").Value = -920571.25
$ws.Range("H61").Value = 2364.0908
$ws.Range("I61").Value = 2364.0908
$ws.Range("K61").Value = 2364.0908
$ws.Range("M61").Value = -2162.0908
$ws.Range("H113").Value = 2364.0908
$ws.Range("I113").Value = 2364.0908
$ws.Range("K113").Value = 2364.0908
$ws.Range("M113").Value = -194.0907999999999
$ws.Range("H132").Value = 6169.8696
$ws.Range("I132").Value = 3557.3845
$ws.Range("J132").Value = 9566.1
$ws.Range("K132").Value = 10672.1535
$ws.Range("L132").Value = 28698.3
$ws.Range("M132").Value = -8142.1535
$ws.Range("N132").Value = -33758.3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 79999.75
$ws.Range("J46").Value = 79999.75
$ws.Range("L46").Value = 79999.75
$ws.Range("N46").Value = -80461.75
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H134").Value = 79999.75
$ws.Range("J134").Value = 79999.75
$ws.Range("L134").Value = 239999.25
$ws.Range("N134").Value = -245069.25
$ws.Range("H136").Value = 6534.6064
$ws.Range("I136").Value = 7448.6206
$ws.Range("K136").Value = 22345.8618
$ws.Range("M136").Value = -19795.8618

Write-Host "Applied all changes"